# Update column F (dSF) values on Sheet1 to reflect repulled data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 2
    3  = -1
    4  = 2
    5  = -4
    7  = -3
    8  = -2
    9  = -2
    10 = -5
    11 = -4
    12 = 1
    13 = 5
    16 = 3
    17 = 6
    18 = -2
    19 = -5
    20 = 5
    21 = -1
    22 = -4
    23 = 1
    24 = 4
    25 = 1
    26 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
